$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before existing row 387 (Extra/Primera records dated 2021-07-09),
# pushing the old rows 387-413 down to become rows 389-415.
$ws.Rows.Item(387).Insert()
$ws.Rows.Item(387).Insert()

# Row 387 - new weekly record
$ws.Cells.Item(387, 1).Value = 9
$ws.Cells.Item(387, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(387, 3).Value = 'Metropolitana'
$ws.Cells.Item(387, 4).Value = 44746
$ws.Cells.Item(387, 5).Value = 13
$ws.Cells.Item(387, 6).Value = 100112013
$ws.Cells.Item(387, 7).Value = 'Alcachofa'
$ws.Cells.Item(387, 8).Value = 'Española'
$ws.Cells.Item(387, 9).Value = 'Extra'
$ws.Cells.Item(387, 10).Value = 25
$ws.Cells.Item(387, 11).Value = 22000
$ws.Cells.Item(387, 12).Value = 22000
$ws.Cells.Item(387, 13).Value = 22000
$ws.Cells.Item(387, 14).Value = '$/caja 25 unidades'
$ws.Cells.Item(387, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(387, 16).Value = 22000
$ws.Cells.Item(387, 17).Value = 1
$ws.Cells.Item(387, 18).Value = 'Hortaliza'

# Row 388 - new weekly record
$ws.Cells.Item(388, 1).Value = 9
$ws.Cells.Item(388, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(388, 3).Value = 'Metropolitana'
$ws.Cells.Item(388, 4).Value = 44746
$ws.Cells.Item(388, 5).Value = 13
$ws.Cells.Item(388, 6).Value = 100112013
$ws.Cells.Item(388, 7).Value = 'Alcachofa'
$ws.Cells.Item(388, 8).Value = 'Española'
$ws.Cells.Item(388, 9).Value = 'Primera'
$ws.Cells.Item(388, 10).Value = 43
$ws.Cells.Item(388, 11).Value = 20000
$ws.Cells.Item(388, 12).Value = 20000
$ws.Cells.Item(388, 13).Value = 20000
$ws.Cells.Item(388, 14).Value = '$/caja 30 unidades'
$ws.Cells.Item(388, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(388, 16).Value = 667
$ws.Cells.Item(388, 17).Value = 30
$ws.Cells.Item(388, 18).Value = 'Hortaliza'
